$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Shared-string text change: "Ready for handoff" -> "In Translation"
#    This string is referenced from:
#      - Overview!E2, Overview!F2
#      - zh-cn!C2
#      - de-de!C2
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2) Narrow the "Status"-column widths from ~17.22 to ~13.41 (character
#    units) on all three sheets:
#      - Overview columns E and F
#      - zh-cn column C
#      - de-de column C
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
